# RobotBOM.xlsx edit:
#  - "9-DOF Board" link (B8) updated to a new Adafruit product URL.
#  - The caster row (row 16) is swapped from "Caster 3/8"" (Pololu #950,
#    $1.99) to "Caster 1"" (Pololu #2691, $3.95), with a new alternate-caster
#    link (Pololu #952) and a currency-style price cell.
#  - Selection cursor moved to F21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 9-DOF Board link change (row 8) ---
$ws.Range("B8").Value = "https://www.adafruit.com/products/1714"

# --- Caster row (row 16): 3/8" -> 1" ---
$ws.Range("A16").Value = "Caster 1"""
$ws.Range("B16").Value = "https://www.pololu.com/product/2691"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3.95
$ws.Range("D16").Style = "Currency"
$ws.Range("E16").Formula = "=D16*C16"
$ws.Range("F16").Value = "Any caster will work, like this one: https://www.pololu.com/product/952"

# --- Move the active selection to F21 ---
$ws.Range("F21").Select()
